$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns.Item(14).ColumnWidth = $ws.Columns.Item(14).ColumnWidth
$ws.Columns.Item(14).NumberFormat = "_(* #,##0.00_);_(* \(#,##0.00\);_(* ""-""??_);_(@_)"
Write-Host "done"
